$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion rates inside the A1 note ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 5.26 = 21052.91 pesos
✅ 21052.91 pesos = 5.24 = 960.49 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$wsHoja1.Range("A1").Value = $newText.TrimEnd()

# --- Sheet "tasas": update the four rate figures (N10, O10, N12, O12) ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 189.998
$wsTasas.Range("O10").Value = 4000.01
$wsTasas.Range("N12").Value = 4017.95
$wsTasas.Range("O12").Value = 183.309
